# OCEPROJECT-4912 - added load 'click' test for R4R results (#155)
#
# Adds a "click" load-test scenario to the R4RResultsLoad sheet: new
# ActionStatus / Filters columns plus three additional rows modelling
# filtered/toolType/researchArea result-page variants, and makes the
# R4RResultsLoad tab the active one when the workbook is reopened.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R4RResultsLoad")

# --- Column B (ContentType labels) for the new scenario rows -------------
$ws.Range("B2").Value = "R4R Results (all)"
$ws.Range("B3").Value = "R4R Results (tools)"
$ws.Range("B4").Value = "R4R Results (areas)"
$ws.Range("B5").Value = "R4R Results (filtered)"

# --- Column A (Path) for the new rows (row 2's path already existed) -----
$ws.Range("A3").Value = "/research/resources/search?from=0&toolTypes=analysis_tools"
$ws.Range("A4").Value = "/research/resources/search?from=0&researchAreas=cancer_omics"
$ws.Range("A5").Value = "/research/resources/search?from=20&toolSubtypes=modeling&toolSubtypes=r_software&toolTypes=analysis_tools"

# --- New headers -----------------------------------------------------------
$ws.Range("C1").Value = "ActionStatus"
$ws.Range("D1").Value = "Filters"

# --- Column C (ActionStatus) ------------------------------------------------
$ws.Range("C2").Value = "r4r_results|view|none|ra=0;tt=0;rt=0;tst=0|1|"
$ws.Range("C3").Value = "r4r_results|view|none|ra=0;tt=1;rt=0;tst=0|1|"
$ws.Range("C4").Value = "r4r_results|view|none|ra=1;tt=0;rt=0;tst=0|1|"
$ws.Range("C5").Value = "r4r_results|view|none|ra=0;tt=1;rt=0;tst=2|2|"

# --- Column D (Filters) -----------------------------------------------------
$ws.Range("D5").Value = "modeling|r_software|analysis_tools"
$ws.Range("D3").Value = "analysis_tools"
$ws.Range("D4").Value = "cancer_omics"
$ws.Range("D2").Value = "none"

# --- Column widths (best-fit approximations for the new columns) -----------
$ws.Columns.Item(1).ColumnWidth = 106
$ws.Columns.Item(2).ColumnWidth = 19.1
$ws.Columns.Item(3).ColumnWidth = 42.1
$ws.Columns.Item(4).ColumnWidth = 33.42

# --- Selection + make this the active/visible tab on reopen -----------------
$ws.Range("A6").Select() | Out-Null
$ws.Activate() | Out-Null
